# hours update and TAR update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new rows of hours (17 & 18). Column A repeats the
# "1/26/2010" date (stored as text, like the other text dates already in
# the sheet), with the same style as row 16. Columns B/C hold the hours
# and the comment.

# Give rows 17-18 col A the same number format/style as row 16 (date style).
$ws.Range("A16").Copy()
$ws.Range("A17:A18").PasteSpecial(-4122)

# Enter the date as a text formula (keeps it from being reinterpreted as a
# serial date number), then collapse it down to a plain value in place so
# no formula or extra styles remain - matching the rest of the sheet.
$ws.Range("A17").Formula = "=""1/26/2010"""
$ws.Range("A17").Copy()
$ws.Range("A17").PasteSpecial(-4163)

$ws.Range("A18").Formula = "=""1/26/2010"""
$ws.Range("A18").Copy()
$ws.Range("A18").PasteSpecial(-4163)

$ws.Range("B17").Value = 2
$ws.Range("C17").Value = "Group Meeting"

$ws.Range("B18").Value = 1
$ws.Range("C18").Value = "Weekly Meeting"

$ws.Range("A19").Select()
